$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style template rows:
#   row 108 = "dark" banding (A-E styles 8,9,9,10,9)
#   row 100 = "light" banding (A-E styles 11,17,17,18,17)
$darkTemplate  = "A108:E108"
$lightTemplate = "A100:E100"

# Row 123
$ws.Range($lightTemplate).Copy()
$ws.Range("A123:E123").PasteSpecial(-4122)
$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = "Review feedback from TRANSCEND:  Address the changes (documentation, new features, improvements, bugs)"
$ws.Cells.Item(123, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(123, 4).Value = 40037
$ws.Cells.Item(123, 5).Value = "In Progress"
$ws.Rows.Item(123).RowHeight = 31

# Row 124
$ws.Range($lightTemplate).Copy()
$ws.Range("A124:E124").PasteSpecial(-4122)
$ws.Cells.Item(124, 1).Value = 123
$ws.Cells.Item(124, 2).Value = "Review feedback from TRANSCEND: Address the performance and time-out issues - collect sample test files"
$ws.Cells.Item(124, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(124, 4).Value = 40037
$ws.Cells.Item(124, 5).Value = "In Progress"
$ws.Rows.Item(124).RowHeight = 31

# Row 125
$ws.Range($darkTemplate).Copy()
$ws.Range("A125:E125").PasteSpecial(-4122)
$ws.Cells.Item(125, 1).Value = 124
$ws.Cells.Item(125, 2).Value = "Installation Guide Updates: Work with Jill and Dev Team to make sure that the installation guide is updated"
$ws.Cells.Item(125, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(125, 4).Value = 40037
$ws.Cells.Item(125, 5).Value = "Complete"
$ws.Rows.Item(125).RowHeight = 31

# Row 126
$ws.Range($darkTemplate).Copy()
$ws.Range("A126:E126").PasteSpecial(-4122)
$ws.Cells.Item(126, 1).Value = 125
$ws.Cells.Item(126, 2).Value = "Let Jill know what kind of links we use for the install documentation"
$ws.Cells.Item(126, 3).Value = "Mike Hunter"
$ws.Cells.Item(126, 4).Value = 40044
$ws.Cells.Item(126, 5).Value = "Obsolete"
$ws.Rows.Item(126).RowHeight = 16

# Row 127
$ws.Range($lightTemplate).Copy()
$ws.Range("A127:E127").PasteSpecial(-4122)
$ws.Cells.Item(127, 1).Value = 126
$ws.Cells.Item(127, 2).Value = "Start providing a quick daily status email for caArray and caIntegrator until 2.5.1 and 1.4.1 go to Production"
$ws.Cells.Item(127, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(127, 4).Value = 40044
$ws.Cells.Item(127, 5).Value = "In Progress"
$ws.Rows.Item(127).RowHeight = 31

# Row 128
$ws.Range($darkTemplate).Copy()
$ws.Range("A128:E128").PasteSpecial(-4122)
$ws.Cells.Item(128, 1).Value = 127
$ws.Cells.Item(128, 2).Value = "Identify the minimal requirements for the test plans for caArray 2.5.1 and caIntegrator 1.4.1 to get into QA."
$ws.Cells.Item(128, 3).Value = "Marina Omelchenko, Sudha, Preston Wood"
$ws.Cells.Item(128, 4).Value = 40058
$ws.Cells.Item(128, 5).Value = "Complete"
$ws.Rows.Item(128).RowHeight = 31

# Row 129
$ws.Range($darkTemplate).Copy()
$ws.Range("A129:E129").PasteSpecial(-4122)
$ws.Cells.Item(129, 1).Value = 128
$ws.Cells.Item(129, 2).Value = "Email the current status of the caArray POAM issues to JJ to ensure we're on track for all outstanding issues."
$ws.Cells.Item(129, 3).Value = "Mike Hunter"
$ws.Cells.Item(129, 4).Value = 40058
$ws.Cells.Item(129, 5).Value = "Obsolete"
$ws.Rows.Item(129).RowHeight = 31

# Row 130
$ws.Range($darkTemplate).Copy()
$ws.Range("A130:E130").PasteSpecial(-4122)
$ws.Cells.Item(130, 1).Value = 129
$ws.Cells.Item(130, 2).Value = "Send the SQL audit log script to Systems and JJ for a monthly cron job."
$ws.Cells.Item(130, 3).Value = "Winston Cheng"
$ws.Cells.Item(130, 4).Value = 40065
$ws.Cells.Item(130, 5).Value = "Complete"
$ws.Rows.Item(130).RowHeight = 16

# Row 131
$ws.Range($darkTemplate).Copy()
$ws.Range("A131:E131").PasteSpecial(-4122)
$ws.Cells.Item(131, 1).Value = 130
$ws.Cells.Item(131, 2).Value = "Determine the optimal records to keep as public for the caIntegrator appscan on STAGE."
$ws.Cells.Item(131, 3).Value = "Mike Hunter and Abe Evans-El"
$ws.Cells.Item(131, 4).Value = 40065
$ws.Cells.Item(131, 5).Value = "Obsolete"
$ws.Rows.Item(131).RowHeight = 31

# Row 132
$ws.Range($darkTemplate).Copy()
$ws.Range("A132:E132").PasteSpecial(-4122)
$ws.Cells.Item(132, 1).Value = 131
$ws.Cells.Item(132, 2).Value = "Provide language and instructions to Jill for handling the 2.5.0 installation instructions."
$ws.Cells.Item(132, 3).Value = "Juli Klemm"
$ws.Cells.Item(132, 4).Value = 40072
$ws.Cells.Item(132, 5).Value = "Complete"
$ws.Rows.Item(132).RowHeight = 31

# Row 112
$ws.Range($darkTemplate).Copy()
$ws.Range("A112:E112").PasteSpecial(-4122)
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "Review caArray audit log capabilities in the 7/30 status meeting."
$ws.Cells.Item(112, 3).Value = "Mike Hunter and Winston Cheng"
$ws.Cells.Item(112, 4).Value = 40009
$ws.Cells.Item(112, 5).Value = "Complete"
$ws.Rows.Item(112).RowHeight = 16

# Row 116
$ws.Range($darkTemplate).Copy()
$ws.Range("A116:E116").PasteSpecial(-4122)
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "After Abe addresses the performance drop for Agilent Data Sets, share Abe's performance page from the wiki with Eve Shalley"
$ws.Cells.Item(116, 3).Value = "Mike Hunter"
$ws.Cells.Item(116, 4).Value = 40016
$ws.Cells.Item(116, 5).Value = "Complete"
$ws.Rows.Item(116).RowHeight = 31

# Row 118
$ws.Range($darkTemplate).Copy()
$ws.Range("A118:E118").PasteSpecial(-4122)
$ws.Cells.Item(118, 1).Value = 117
$ws.Cells.Item(118, 2).Value = "Follow up with Ann Wiley to find out what she needs from us to update the product landing page."
$ws.Cells.Item(118, 3).Value = "Mike Hunter"
$ws.Cells.Item(118, 4).Value = 40023
$ws.Cells.Item(118, 5).Value = "Complete"
$ws.Rows.Item(118).RowHeight = 31

# Row 119
$ws.Range($darkTemplate).Copy()
$ws.Range("A119:E119").PasteSpecial(-4122)
$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = "Meet to refine the audit log requirements. (Thursday: after 2:00)"
$ws.Cells.Item(119, 3).Value = "Ulli Wagner, Mike Hunter, JJ Pan, Juli Klemm, Winston Cheng"
$ws.Cells.Item(119, 4).Value = 40023
$ws.Cells.Item(119, 5).Value = "Complete"
$ws.Rows.Item(119).RowHeight = 16

# Row 120
$ws.Range($darkTemplate).Copy()
$ws.Range("A120:E120").PasteSpecial(-4122)
$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = "Change Management Request #59 (Java 1.7) upgrade. Review and Respond."
$ws.Cells.Item(120, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(120, 4).Value = 40023
$ws.Cells.Item(120, 5).Value = "Complete"
$ws.Rows.Item(120).RowHeight = 16

# Row 121
$ws.Range($darkTemplate).Copy()
$ws.Range("A121:E121").PasteSpecial(-4122)
$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = "Grid dependencies for caArray and caIntegrator. Review and Respond once Ulli has the official request."
$ws.Cells.Item(121, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(121, 4).Value = 40030
$ws.Cells.Item(121, 5).Value = "Complete"
$ws.Rows.Item(121).RowHeight = 31

# Row 122
$ws.Range($darkTemplate).Copy()
$ws.Range("A122:E122").PasteSpecial(-4122)
$ws.Cells.Item(122, 1).Value = 121
$ws.Cells.Item(122, 2).Value = "Status update on the caArray QA tier. Jacob will forward  Shine, his latest communication with Winston."
$ws.Cells.Item(122, 3).Value = "Mike Hunter and Shine Jacob"
$ws.Cells.Item(122, 4).Value = 40030
$ws.Cells.Item(122, 5).Value = "Complete"
$ws.Rows.Item(122).RowHeight = 31

# Row 133
$ws.Range($lightTemplate).Copy()
$ws.Range("A133:E133").PasteSpecial(-4122)
$ws.Cells.Item(133, 1).Value = 132
$ws.Cells.Item(133, 2).Value = ""
$ws.Cells.Item(133, 3).Value = ""
$ws.Cells.Item(133, 4).Value = ""
$ws.Cells.Item(133, 5).Value = ""
$ws.Rows.Item(133).RowHeight = 16

# Row 134
$ws.Range($lightTemplate).Copy()
$ws.Range("A134:E134").PasteSpecial(-4122)
$ws.Cells.Item(134, 1).Value = 133
$ws.Cells.Item(134, 2).Value = ""
$ws.Cells.Item(134, 3).Value = ""
$ws.Cells.Item(134, 4).Value = ""
$ws.Cells.Item(134, 5).Value = ""
$ws.Rows.Item(134).RowHeight = 16

# Row 135
$ws.Range($lightTemplate).Copy()
$ws.Range("A135:E135").PasteSpecial(-4122)
$ws.Cells.Item(135, 1).Value = 134
$ws.Cells.Item(135, 2).Value = ""
$ws.Cells.Item(135, 3).Value = ""
$ws.Cells.Item(135, 4).Value = ""
$ws.Cells.Item(135, 5).Value = ""
$ws.Rows.Item(135).RowHeight = 16

# Row 136
$ws.Range($lightTemplate).Copy()
$ws.Range("A136:E136").PasteSpecial(-4122)
$ws.Cells.Item(136, 1).Value = 135
$ws.Cells.Item(136, 2).Value = ""
$ws.Cells.Item(136, 3).Value = ""
$ws.Cells.Item(136, 4).Value = ""
$ws.Cells.Item(136, 5).Value = ""
$ws.Rows.Item(136).RowHeight = 16

# Row 137
$ws.Range($lightTemplate).Copy()
$ws.Range("A137:E137").PasteSpecial(-4122)
$ws.Cells.Item(137, 1).Value = 136
$ws.Cells.Item(137, 2).Value = ""
$ws.Cells.Item(137, 3).Value = ""
$ws.Cells.Item(137, 4).Value = ""
$ws.Cells.Item(137, 5).Value = ""
$ws.Rows.Item(137).RowHeight = 16

# Row 138
$ws.Range($lightTemplate).Copy()
$ws.Range("A138:E138").PasteSpecial(-4122)
$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 2).Value = ""
$ws.Cells.Item(138, 3).Value = ""
$ws.Cells.Item(138, 4).Value = ""
$ws.Cells.Item(138, 5).Value = ""
$ws.Rows.Item(138).RowHeight = 16

# Row 139
$ws.Range($lightTemplate).Copy()
$ws.Range("A139:E139").PasteSpecial(-4122)
$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = ""
$ws.Cells.Item(139, 3).Value = ""
$ws.Cells.Item(139, 4).Value = ""
$ws.Cells.Item(139, 5).Value = ""
$ws.Rows.Item(139).RowHeight = 16

# Row 140
$ws.Range($lightTemplate).Copy()
$ws.Range("A140:E140").PasteSpecial(-4122)
$ws.Cells.Item(140, 1).Value = 139
$ws.Cells.Item(140, 2).Value = ""
$ws.Cells.Item(140, 3).Value = ""
$ws.Cells.Item(140, 4).Value = ""
$ws.Cells.Item(140, 5).Value = ""
$ws.Rows.Item(140).RowHeight = 16

# Row 141
$ws.Range($lightTemplate).Copy()
$ws.Range("A141:E141").PasteSpecial(-4122)
$ws.Cells.Item(141, 1).Value = 140
$ws.Cells.Item(141, 2).Value = ""
$ws.Cells.Item(141, 3).Value = ""
$ws.Cells.Item(141, 4).Value = ""
$ws.Cells.Item(141, 5).Value = ""
$ws.Rows.Item(141).RowHeight = 16

# Row 142
$ws.Range($lightTemplate).Copy()
$ws.Range("A142:E142").PasteSpecial(-4122)
$ws.Cells.Item(142, 1).Value = 141
$ws.Cells.Item(142, 2).Value = ""
$ws.Cells.Item(142, 3).Value = ""
$ws.Cells.Item(142, 4).Value = ""
$ws.Cells.Item(142, 5).Value = ""
$ws.Rows.Item(142).RowHeight = 16

# Row 143
$ws.Range($lightTemplate).Copy()
$ws.Range("A143:E143").PasteSpecial(-4122)
$ws.Cells.Item(143, 1).Value = 142
$ws.Cells.Item(143, 2).Value = ""
$ws.Cells.Item(143, 3).Value = ""
$ws.Cells.Item(143, 4).Value = ""
$ws.Cells.Item(143, 5).Value = ""
$ws.Rows.Item(143).RowHeight = 16

# Row 144
$ws.Range($lightTemplate).Copy()
$ws.Range("A144:E144").PasteSpecial(-4122)
$ws.Cells.Item(144, 1).Value = 143
$ws.Cells.Item(144, 2).Value = ""
$ws.Cells.Item(144, 3).Value = ""
$ws.Cells.Item(144, 4).Value = ""
$ws.Cells.Item(144, 5).Value = ""
$ws.Rows.Item(144).RowHeight = 16

# Row 145
$ws.Range($lightTemplate).Copy()
$ws.Range("A145:E145").PasteSpecial(-4122)
$ws.Cells.Item(145, 1).Value = 144
$ws.Cells.Item(145, 2).Value = ""
$ws.Cells.Item(145, 3).Value = ""
$ws.Cells.Item(145, 4).Value = ""
$ws.Cells.Item(145, 5).Value = ""
$ws.Rows.Item(145).RowHeight = 16

# Row 146
$ws.Range($lightTemplate).Copy()
$ws.Range("A146:E146").PasteSpecial(-4122)
$ws.Cells.Item(146, 1).Value = 145
$ws.Cells.Item(146, 2).Value = ""
$ws.Cells.Item(146, 3).Value = ""
$ws.Cells.Item(146, 4).Value = ""
$ws.Cells.Item(146, 5).Value = ""
$ws.Rows.Item(146).RowHeight = 16

# Row 147
$ws.Range($lightTemplate).Copy()
$ws.Range("A147:E147").PasteSpecial(-4122)
$ws.Cells.Item(147, 1).Value = 146
$ws.Cells.Item(147, 2).Value = ""
$ws.Cells.Item(147, 3).Value = ""
$ws.Cells.Item(147, 4).Value = ""
$ws.Cells.Item(147, 5).Value = ""
$ws.Rows.Item(147).RowHeight = 16

# Row 148
$ws.Range($lightTemplate).Copy()
$ws.Range("A148:E148").PasteSpecial(-4122)
$ws.Cells.Item(148, 1).Value = 147
$ws.Cells.Item(148, 2).Value = ""
$ws.Cells.Item(148, 3).Value = ""
$ws.Cells.Item(148, 4).Value = ""
$ws.Cells.Item(148, 5).Value = ""
$ws.Rows.Item(148).RowHeight = 16

# Row 149
$ws.Range($lightTemplate).Copy()
$ws.Range("A149:E149").PasteSpecial(-4122)
$ws.Cells.Item(149, 1).Value = 148
$ws.Cells.Item(149, 2).Value = ""
$ws.Cells.Item(149, 3).Value = ""
$ws.Cells.Item(149, 4).Value = ""
$ws.Cells.Item(149, 5).Value = ""
$ws.Rows.Item(149).RowHeight = 16

# Row 150
$ws.Range($lightTemplate).Copy()
$ws.Range("A150:E150").PasteSpecial(-4122)
$ws.Cells.Item(150, 1).Value = 149
$ws.Cells.Item(150, 2).Value = ""
$ws.Cells.Item(150, 3).Value = ""
$ws.Cells.Item(150, 4).Value = ""
$ws.Cells.Item(150, 5).Value = ""
$ws.Rows.Item(150).RowHeight = 16

# Row 151
$ws.Range($lightTemplate).Copy()
$ws.Range("A151:E151").PasteSpecial(-4122)
$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = ""
$ws.Cells.Item(151, 3).Value = ""
$ws.Cells.Item(151, 4).Value = ""
$ws.Cells.Item(151, 5).Value = ""
$ws.Rows.Item(151).RowHeight = 16

# Update the active view: clear the pinned top-left cell and move the selection
# to the newly-added rows (matches the refreshed status-meeting viewport).
$ws.Range("A128:E132").Select()
